# Add random forest classifier model row + a new "Model" column with hyperlinked notebook names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new blank column C ("Model") ---------------------------------
$ws.Columns("C").Insert()

# --- 2. New header for the inserted column ------------------------------------
$ws.Range("C1").Value = "Model"

# --- 3. Fill in row 4 with the new Random Forest trial -------------------------
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "bank_marketing_random_forest"
$ws.Range("C4").Value = "Random Forest Classifier"
$ws.Range("D4").Value = "N/A"
$ws.Range("E4").Value = "N/A"
$ws.Range("F4").Value = "N/A"
$ws.Range("G4").Value = "N/A"
$ws.Range("H4").Value = 0.91

# --- 4. Hyperlink the three notebook-name cells --------------------------------
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/dakog/project4/blob/main/bank_marketing.ipynb")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/dakog/project4/blob/main/bank_marketing_auto_opt1.ipynb")
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/dakog/project4/blob/main/bank_marketing_random_forest.ipynb")

# --- 5. Cosmetic view/format refresh -------------------------------------------
[void]$ws.Range("A1:H7").EntireColumn.AutoFit()

[void]$ws.Range("J7").Select()
